# Auto-generated edit script: updates cryptos price/volume columns (D, E)
# to match the refreshed coinranking.com snapshot referenced in the commit message.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.792.76'
$ws.Range('E2').Value = '  +1.05%  '
$ws.Range('D3').Value = '2.492.13'
$ws.Range('E3').Value = '  +0.40%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Formula = "'586.64"
$ws.Range('E5').Value = '  +0.31%  '
$ws.Range('D6').Formula = "'176.95"
$ws.Range('E6').Value = '  +3.33%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('E8').Value = '  +0.31%  '
$ws.Range('E9').Value = '  +3.88%  '
$ws.Range('E10').Value = '  +0.16%  '
$ws.Range('E11').Value = '  +2.56%  '
$ws.Range('E12').Value = '  +0.20%  '
$ws.Range('D13').Value = '2.948.20'
$ws.Range('E13').Value = '  +0.48%  '
$ws.Range('D14').Formula = "'25.68"
$ws.Range('E14').Value = '  +1.20%  '
$ws.Range('D15').Value = '67.717.56'
$ws.Range('E15').Value = '  +1.14%  '
$ws.Range('E16').Value = '  +0.91%  '
$ws.Range('D17').Value = '2.485.47'
$ws.Range('E17').Value = '  +0.15%  '
$ws.Range('D18').Formula = "'7.52"
$ws.Range('E18').Value = '  +1.74%  '
$ws.Range('D19').Formula = "'10.97"
$ws.Range('E19').Value = '  +0.06%  '
$ws.Range('D20').Formula = "'350.43"
$ws.Range('E20').Value = '  +0.16%  '
$ws.Range('D21').Formula = "'4.10"
$ws.Range('E21').Value = '  +2.30%  '
$ws.Range('D23').Formula = "'70.76"
$ws.Range('E23').Value = '  +3.30%  '
$ws.Range('E24').Value = '  +0.95%  '
$ws.Range('D25').Formula = "'1.74"
$ws.Range('E25').Value = '  -2.46%  '
$ws.Range('D26').Formula = "'9.11"
$ws.Range('E26').Value = '  -1.28%  '
$ws.Range('E27').Value = '  +0.25%  '
$ws.Range('D28').Formula = "'0.996"
$ws.Range('E28').Value = '  -0.29%  '
$ws.Range('D29').Value = '0.0₃0903'
$ws.Range('E29').Value = '  +0.58%  '
$ws.Range('D30').Formula = "'505.15"
$ws.Range('E30').Value = '  -1.00%  '
$ws.Range('D31').Formula = "'7.83"
$ws.Range('E31').Value = '  +1.78%  '
$ws.Range('E32').Value = '  +2.60%  '
$ws.Range('E33').Value = '  +0.64%  '
$ws.Range('E34').Value = '  +0.02%  '
$ws.Range('E35').Value = '  +3.70%  '
$ws.Range('D36').Formula = "'162.87"
$ws.Range('E36').Value = '  +1.78%  '
$ws.Range('D37').Formula = "'18.68"
$ws.Range('E37').Value = '  -0.08%  '
$ws.Range('D38').Formula = "'18.31"
$ws.Range('E38').Value = '  +0.56%  '
$ws.Range('E39').Value = '  +0.46%  '
$ws.Range('E40').Value = '  +0.02%  '
$ws.Range('E41').Value = '  +3.72%  '
$ws.Range('E42').Value = '  +0.36%  '
$ws.Range('E43').Value = '  +0.98%  '
$ws.Range('D44').Formula = "'2.42"
$ws.Range('E44').Value = '  +1.89%  '
$ws.Range('D45').Formula = "'144.66"
$ws.Range('E45').Value = '  +1.32%  '
$ws.Range('E46').Value = '  +2.34%  '
$ws.Range('D47').Formula = "'0.514"
$ws.Range('E47').Value = '  +0.05%  '
$ws.Range('E48').Value = '  +1.48%  '
$ws.Range('D49').Formula = "'0.0742"
$ws.Range('E49').Value = '  +1.81%  '
$ws.Range('E50').Value = '  +1.28%  '
$ws.Range('E51').Value = '  +0.64%  '
